$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the first data row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-02 01:17:19"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first data row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-02 01:17:14"
$wsZhCn.Range("K2").Value = "2016-09-02 01:17:32"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first data row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-02 01:17:19"
$wsDeDe.Range("K2").Value = "2016-09-02 01:17:39"
